$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the main input probability table (C3:F6) ---
$ws.Range("C3").Value = 0.73
$ws.Range("D4").Value = 0.25
$ws.Range("C6").Value = 0.18
$ws.Range("D6").Value = 0

# --- Update the "Credibility" indicator tables (left: J12:M15, right: P12:S15) ---
$ws.Range("J13").Value = 1
$ws.Range("K15").Value = 0
$ws.Range("Q15").Value = 0

# --- Update the active cell selection to match the author's last position ---
$ws.Range("K16").Select()
